# livehta_2301_Data.xlsx — "Non-Oncology Excel reports comparison"
#
# The source diff shows two data cells on Sheet1 being repointed to a
# (slightly re-worded, hyphen-no-space) "CompleteExcelReport" filename
# string, plus the resulting sheet-view selection landing on the second
# of those two cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newReportName = "CompleteExcelReport-LIVEHTA Automation-Test_NonOncology_Automation_3-Clinical-2023_"

# Column I ("ExpectedFilenames") on rows 3 and 11 both get the new value.
$ws.Range("I3").Value = $newReportName
$ws.Range("I11").Value = $newReportName

# Leave the selection on I11, matching the saved view state.
[void]$ws.Range("I11").Select()
